# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" colours (currently wired to the
#                            notes master)
#   ppt/theme/theme2.xml -> "Integral" colours (currently wired to the slide
#                            master / presentation)
#
# The authored change swaps which theme the slide master uses for its
# 12-colour scheme (and, symmetrically, what the notes master uses) -- i.e.
# the deck's visible/working theme flips from the "Integral" palette to the
# stock "Office Theme" palette. Reproduce that by rewriting every slot of
# the slide master's ThemeColorScheme to the Office Theme RGB values.

function Convert-HexToComRgb {
    param([string]$Hex)
    $r = [Convert]::ToInt32($Hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($Hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($Hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Office Theme colour scheme, in the standard dk1/lt1/dk2/lt2/accent1-6/
# hlink/folHlink order used by ThemeColorScheme.Item(1..12).
$officeThemeHex = @(
    "000000", # 1  dk1
    "FFFFFF", # 2  lt1
    "44546A", # 3  dk2
    "E7E6E6", # 4  lt2
    "5B9BD5", # 5  accent1
    "ED7D31", # 6  accent2
    "A5A5A5", # 7  accent3
    "FFC000", # 8  accent4
    "4472C4", # 9  accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hlink
    "954F72"  # 12 folHlink
)

$tcs = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeThemeHex.Count; $i++) {
    $tcs.Item($i).RGB = Convert-HexToComRgb $officeThemeHex[$i - 1]
}
